$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "0.691") need to be
# forced to Text format first, otherwise Excel auto-converts them to numbers
# (these columns hold pre-formatted display strings, same as the source data).
$textForceCells = @("D5", "D6", "D8", "D11", "D13", "D15", "D19", "D21", "D22", "D24", "D27", "D29", "D31", "D37", "D38", "D39", "D41", "D42", "D44", "D47")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin list values (price / 1h volume refresh, plus the
# BNB/XRP and InjectiveProtocol/Aave row-order swaps).
$ws.Range("D2").Value = '35.370.80'
$ws.Range("E2").Value = '  +1.87%  '
$ws.Range("D3").Value = '1.883.02'
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("B5").Value = 'XRP'
$ws.Range("C5").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D5").Value = '0.691'
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '245.67'
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '43.21'
$ws.Range("E8").Value = '  +2.87%  '
$ws.Range("E9").Value = '  +2.81%  '
$ws.Range("E10").Value = '  +7.54%  '
$ws.Range("D11").Value = '0.0743'
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("E12").Value = '  +1.71%  '
$ws.Range("D13").Value = '13.77'
$ws.Range("E13").Value = '  +8.04%  '
$ws.Range("D14").Value = '2.156.52'
$ws.Range("E14").Value = '  +0.68%  '
$ws.Range("D15").Value = '0.770'
$ws.Range("E15").Value = '  +8.42%  '
$ws.Range("D17").Value = '1.891.01'
$ws.Range("E17").Value = '  +1.02%  '
$ws.Range("D18").Value = '35.338.23'
$ws.Range("E18").Value = '  +1.88%  '
$ws.Range("D19").Value = '73.40'
$ws.Range("E19").Value = '  +1.21%  '
$ws.Range("E20").Value = '  +1.28%  '
$ws.Range("D21").Value = '244.39'
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").Value = '12.81'
$ws.Range("E22").Value = '  +1.38%  '
$ws.Range("E23").Value = '  +4.89%  '
$ws.Range("D24").Value = '2.64'
$ws.Range("E24").Value = '  +9.01%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  -3.18%  '
$ws.Range("D27").Value = '165.58'
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("E28").Value = '  +3.28%  '
$ws.Range("D29").Value = '18.26'
$ws.Range("E29").Value = '  +0.71%  '
$ws.Range("E30").Value = '  +0.43%  '
$ws.Range("D31").Value = '0.0595'
$ws.Range("E31").Value = '  +3.77%  '
$ws.Range("E32").Value = '  +0.23%  '
$ws.Range("E33").Value = '  +20.18%  '
$ws.Range("E34").Value = '  +0.81%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  -12.86%  '
$ws.Range("D37").Value = '0.853'
$ws.Range("E37").Value = '  +3.61%  '
$ws.Range("D38").Value = '1.94'
$ws.Range("E38").Value = '  -1.70%  '
$ws.Range("D39").Value = '0.0718'
$ws.Range("E39").Value = '  +9.09%  '
$ws.Range("E40").Value = '  +5.50%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '97.69'
$ws.Range("E41").Value = '  +0.33%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").Value = '17.10'
$ws.Range("E42").Value = '  +1.72%  '
$ws.Range("E43").Value = '  -0.35%  '
$ws.Range("D44").Value = '13.72'
$ws.Range("E44").Value = '  +13.34%  '
$ws.Range("D45").Value = '1.320.99'
$ws.Range("E45").Value = '  +3.19%  '
$ws.Range("E46").Value = '  +2.96%  '
$ws.Range("D47").Value = '0.0809'
$ws.Range("E47").Value = '  +2.91%  '
$ws.Range("E48").Value = '  +0.09%  '
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("E50").Value = '  -2.49%  '
$ws.Range("D51").Value = '2.056.62'
$ws.Range("E51").Value = '  +0.38%  '
